# Applies the 10.2.1.xlsx update:
#  - fills in column J (year 2021) values for rows 19-26, copying the
#    number format / font already used by column I in that block so the
#    new cells render like their neighbours
#  - fills in the previously-blank total cell J27
#  - moves the saved selection from L27 to N8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    19 = 12.434613462352335
    20 = 16.80050595536094
    21 = 11.282963378125267
    22 = 25.042808754677555
    23 = 3.2011163356916352
    24 = 13.523574517571838
    25 = 6.1196997869329204
    26 = 5.9488136666578013
}

foreach ($row in $newValues.Keys) {
    $srcCell = $ws.Range("I$row")
    $dstCell = $ws.Range("J$row")
    $srcCell.Copy($dstCell)
    $dstCell.Value = $newValues[$row]
}

$ws.Range("J27").Value = 5.2451982064110645

$ws.Range("N8").Select() | Out-Null
